$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the last row (row 28 -> 2025Q2) with refreshed recurrence metrics
$ws.Range("C28").Value = 189
$ws.Range("D28").Value = 26
$ws.Range("E28").Value = 163
$ws.Range("F28").Value = 4.049844236760125
